$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like pure numbers (e.g. "559.11") would
# otherwise be auto-converted to a numeric cell type by Excel. Force
# each such cell to remain text (matching the original inlineStr/text
# cell type in the workbook) by temporarily applying a Text number
# format right before assigning the value, then clearing the format
# again afterwards so no stray style index is left behind.

$ws.Range("D2").Value = "65.341.90"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.387.09"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.68"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.375.91"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.38"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000278"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "3.932.60"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "3.378.81"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "65.402.07"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.77"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.90"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.26"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.91"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.71"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("E30").Value = "  +3.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.54"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.47"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.91"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "573.78"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.59"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.14%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.88"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.373"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").Value = "0.0₃0738"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").Value = "3.129.47"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.79"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0416"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.133"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.16"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.06"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.40"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.75%  "
